# Updated cryptos list on Sat Mar 23 11:12:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

# Prevent Excel from auto-converting numeric-looking price strings to real numbers
# (target keeps these as literal text, matching the original inline-string cells).
# NOTE: a multi-area union Range("D5,D6,...").NumberFormat=... only applies to the
# first area on this host, so each address is set individually in a loop.
$textCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D14", "D16", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D39", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = '64.347.99'
$ws.Range("E2").Value = '  -0.74%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.352.46'
$ws.Range("E3").Value = '  -2.45%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.09%  '

# Row 5 - BNB
$ws.Range("D5").Value = '555.74'
$ws.Range("E5").Value = '  -2.75%  '

# Row 6 - Solana
$ws.Range("D6").Value = '175.26'
$ws.Range("E6").Value = '  +0.28%  '

# Row 7 - XRP
$ws.Range("E7").Value = '  -0.54%  '

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = '3.342.97'
$ws.Range("E8").Value = '  -2.57%  '

# Row 10 - Cardano->Dogecoin
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  +3.28%  '

# Row 11 - Dogecoin->Cardano
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '0.627'
$ws.Range("E11").Value = '  +0.73%  '

# Row 12 - Avalanche
$ws.Range("D12").Value = '54.49'
$ws.Range("E12").Value = '  -0.46%  '

# Row 13 - ShibaInu
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").Value = '  +0.48%  '

# Row 14 - Polkadot
$ws.Range("D14").Value = '9.08'
$ws.Range("E14").Value = '  -0.19%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '3.885.92'
$ws.Range("E15").Value = '  -2.43%  '

# Row 16 - Chainlink
$ws.Range("D16").Value = '18.37'
$ws.Range("E16").Value = '  +1.82%  '

# Row 17 - TRON
$ws.Range("E17").Value = '  -1.86%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '3.347.53'
$ws.Range("E18").Value = '  -2.55%  '

# Row 19 - Uniswap
$ws.Range("D19").Value = '11.84'
$ws.Range("E19").Value = '  +0.07%  '

# Row 20 - WrappedBTC
$ws.Range("D20").Value = '64.269.40'
$ws.Range("E20").Value = '  -0.78%  '

# Row 21 - Polygon
$ws.Range("D21").Value = '0.983'
$ws.Range("E21").Value = '  -0.21%  '

# Row 22 - BitcoinCash
$ws.Range("D22").Value = '455.77'
$ws.Range("E22").Value = '  +12.05%  '

# Row 23 - Toncoin
$ws.Range("D23").Value = '4.86'
$ws.Range("E23").Value = '  +9.87%  '

# Row 24 - PancakeSwap
$ws.Range("D24").Value = '4.08'
$ws.Range("E24").Value = '  -2.56%  '

# Row 25 - Litecoin
$ws.Range("D25").Value = '85.82'
$ws.Range("E25").Value = '  +2.86%  '

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = '13.35'
$ws.Range("E26").Value = '  +0.75%  '

# Row 27 - RenderToken
$ws.Range("D27").Value = '10.98'
$ws.Range("E27").Value = '  +1.78%  '

# Row 28 - ImmutableX
$ws.Range("D28").Value = '2.85'
$ws.Range("E28").Value = '  +1.81%  '

# Row 29 - Filecoin
$ws.Range("D29").Value = '8.78'
$ws.Range("E29").Value = '  -1.74%  '

# Row 30 - EthereumClassic
$ws.Range("D30").Value = '30.01'
$ws.Range("E30").Value = '  +0.63%  '

# Row 31 - NEARProtocol
$ws.Range("D31").Value = '6.67'
$ws.Range("E31").Value = '  +1.78%  '

# Row 32 - Bittensor
$ws.Range("D32").Value = '585.44'
$ws.Range("E32").Value = '  +0.17%  '

# Row 33 - Cosmos
$ws.Range("D33").Value = '11.48'
$ws.Range("E33").Value = '  -0.29%  '

# Row 34 - Hedera
$ws.Range("E34").Value = '  -0.15%  '

# Row 35 - OKB
$ws.Range("D35").Value = '58.69'
$ws.Range("E35").Value = '  -1.50%  '

# Row 36 - Dai
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.09%  '

# Row 37 - Kaspa
$ws.Range("E37").Value = '  -7.84%  '

# Row 38 - Stacks
$ws.Range("D38").Value = '3.52'
$ws.Range("E38").Value = '  -0.57%  '

# Row 39 - InjectiveProtocol
$ws.Range("D39").Value = '35.72'
$ws.Range("E39").Value = '  -1.01%  '

# Row 40 - PEPE
$ws.Range("D40").Value = "0.0{0}0755" -f $sub3
$ws.Range("E40").Value = '  -1.59%  '

# Row 41 - TheGraph
$ws.Range("D41").Value = '0.375'
$ws.Range("E41").Value = '  -0.21%  '

# Row 42 - Maker
$ws.Range("D42").Value = '3.094.18'
$ws.Range("E42").Value = '  -2.62%  '

# Row 43 - FirstDigitalUSD
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.14%  '

# Row 44 - ThetaToken->Fetch.AI
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.54'
$ws.Range("E44").Value = '  +1.57%  '

# Row 45 - Fetch.AI->ThetaToken
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").Value = '2.79'
$ws.Range("E45").Value = '  -4.50%  '

# Row 46 - ApeXProtocol->VeChain
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").Value = '  +0.53%  '

# Row 47 - VeChain->ApeXProtocol
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.18'
$ws.Range("E47").Value = '  -2.11%  '

# Row 48 - Stellar
$ws.Range("D48").Value = '0.131'
$ws.Range("E48").Value = '  +0.19%  '

# Row 49 - WEMIXToken
$ws.Range("D49").Value = '2.58'
$ws.Range("E49").Value = '  -1.81%  '

# Row 50 - THORChain
$ws.Range("D50").Value = '8.35'
$ws.Range("E50").Value = '  -0.78%  '

# Row 51 - Monero
$ws.Range("D51").Value = '135.35'
$ws.Range("E51").Value = '  -1.29%  '
